$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the affected cells stay as plain text (they hold values like
# "641,530,687,575" or "0.76" which must NOT be reinterpreted as numbers).

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "image_20250807110238_ppp0.jpg"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "641,530,687,575"
$ws.Range("J16").NumberFormat = "@"
$ws.Range("J16").Value = "0.76"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "image_20250807110238_ppp0.jpg"
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "793,481,831,527"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "image_20250808221835_ppp0.jpg"
$ws.Range("I18").NumberFormat = "@"
$ws.Range("I18").Value = "1182,405,1231,455"
$ws.Range("J18").NumberFormat = "@"
$ws.Range("J18").Value = "0.76"
